# Weekly update: insert a new price-report row for
# "Feria Lagunitas de Puerto Montt - Ciboulette" dated 2023-06-12,
# pushing the existing rows 242:370 down to 243:371.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 242 (Excel shifts 242:370 -> 243:371
# and copies the formatting of the row above, which already carries the
# date-cell style used throughout column D).
$ws.Range("A242").EntireRow.Insert()

# Populate the newly inserted row with this week's values.
$ws.Range("A242").Value = 4
$ws.Range("B242").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C242").Value = "Los Lagos"
$ws.Range("D242").Value = [DateTime]"2023-06-12"
$ws.Range("E242").Value = 10
$ws.Range("F242").Value = 100112039
$ws.Range("G242").Value = "Ciboulette"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 80
$ws.Range("K242").Value = 2500
$ws.Range("L242").Value = 2500
$ws.Range("M242").Value = 2500
$ws.Range("N242").Value = "$/docena de atados"
$ws.Range("O242").Value = "Región Metropolitana"
$ws.Range("P242").Value = 833
$ws.Range("Q242").Value = 3
$ws.Range("R242").Value = "Hortaliza"
